# WRI updates for RTMF
# Updates the Recipient Transportation Mode Fractions workbook:
#  - RTMF-passengers: row 7 (ships) now shifts 0.33 to aircraft and 0.33 to rail
#  - RTMF-freight: row 2 (LDVs) now shifts 0.5 to aircraft and 0.5 to rail
#  - RTMF-freight: row 3 (HDVs) rail/ships shift changes from 0.9/0.1 to 0.8/0.2
#  - RTMF-freight: row 7 (ships) now shifts 0.33 to aircraft and 0.33 to rail
#  - Leaves the workbook with the "RTMF-passengers" sheet active (matching
#    the saved view state) and records the last selected cell on each sheet.

$wb = $excel.ActiveWorkbook

$wsAbout      = $wb.Worksheets.Item("About")
$wsPassengers = $wb.Worksheets.Item("RTMF-passengers")
$wsFreight    = $wb.Worksheets.Item("RTMF-freight")

# --- RTMF-passengers: row 7 (ships row) -----------------------------------
$wsPassengers.Range("C7").Value = 0.33
$wsPassengers.Range("E7").Value = 0.33

# --- RTMF-freight: row 2 (LDVs row) ----------------------------------------
$wsFreight.Range("C2").Value = 0.5
$wsFreight.Range("E2").Value = 0.5

# --- RTMF-freight: row 3 (HDVs row) ----------------------------------------
$wsFreight.Range("E3").Value = 0.8
$wsFreight.Range("F3").Value = 0.2

# --- RTMF-freight: row 7 (ships row) ----------------------------------------
$wsFreight.Range("C7").Value = 0.33
$wsFreight.Range("E7").Value = 0.33

# --- View / selection state --------------------------------------------------
# "RTMF-freight" keeps its own last-selected cell too.
$wsFreight.Activate()
$wsFreight.Range("F7").Select()

# "About" scrolled down with C25 as the last selected cell.
$wsAbout.Activate()
$wsAbout.Range("C25").Select()

# "RTMF-passengers" ends up the active/selected tab with F7 selected.
$wsPassengers.Activate()
$wsPassengers.Range("F7").Select()
